$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 870.8570999999999
$ws.Range("I18").Value = 870.8570999999999
$ws.Range("K18").Value = 870.8570999999999
$ws.Range("M18").Value = -586.8570999999999
$ws.Range("H96").Value = 732.25
$ws.Range("I96").Value = 735.75
$ws.Range("J96").Value = 721.75
$ws.Range("K96").Value = 2207.25
$ws.Range("L96").Value = 2165.25
$ws.Range("M96").Value = -834.25
$ws.Range("N96").Value = -4911.25
$ws.Range("H99").Value = 385.7143
$ws.Range("J99").Value = 600.2
$ws.Range("L99").Value = 1800.6
$ws.Range("N99").Value = -4796.6
$ws.Range("H100").Value = 1898.45
$ws.Range("I100").Value = 606.8182
$ws.Range("J100").Value = 3477.111
$ws.Range("K100").Value = 606.8182
$ws.Range("L100").Value = 3477.111
$ws.Range("M100").Value = -65.81820000000005
$ws.Range("N100").Value = -4559.111
$ws.Range("H120").Value = 74992.5
$ws.Range("J120").Value = 74992.5
$ws.Range("L120").Value = 74992.5
$ws.Range("N120").Value = -84668.5
$ws.Range("H123").Value = 86165.71000000001
$ws.Range("J123").Value = 86165.71000000001
$ws.Range("L123").Value = 86165.71000000001
$ws.Range("N123").Value = -95965.71000000001
$ws.Range("H131").Value = 1061.6666
$ws.Range("I131").Value = 1061.6666
$ws.Range("K131").Value = 3184.9998
$ws.Range("M131").Value = 1855.0002
$ws.Range("H132").Value = 1846.1428
$ws.Range("I132").Value = 800.28125
$ws.Range("K132").Value = 2400.84375
$ws.Range("M132").Value = 129.15625
$ws.Range("H133").Value = 92928.57000000001
$ws.Range("J133").Value = 92928.57000000001
$ws.Range("L133").Value = 92928.57000000001
$ws.Range("N133").Value = -103048.57
$ws.Range("H134").Value = 61781.816
$ws.Range("J134").Value = 65869.44500000001
$ws.Range("L134").Value = 65869.44500000001
$ws.Range("N134").Value = -76009.44500000001
$ws.Range("H136").Value = 72857.2
$ws.Range("J136").Value = 82071.5
$ws.Range("L136").Value = 82071.5
$ws.Range("N136").Value = -92271.5
$ws.Range("H138").Value = 2281.0933
$ws.Range("I138").Value = 1565.037
$ws.Range("K138").Value = 4695.111
$ws.Range("M138").Value = 444.8890000000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7696.1113
$ws.Range("I45").Value = 8235.4
$ws.Range("K45").Value = 8235.4
$ws.Range("M45").Value = -7858.4
$ws.Range("H74").Value = 35634
$ws.Range("I74").Value = 44957.22
$ws.Range("K74").Value = 44957.22
$ws.Range("M74").Value = -44083.22
$ws.Range("H77").Value = 35634
$ws.Range("I77").Value = 44957.22
$ws.Range("K77").Value = 224786.1
$ws.Range("M77").Value = -220418.1
$ws.Range("H108").Value = 89992.5
$ws.Range("J108").Value = 89992.5
$ws.Range("L108").Value = 89992.5
$ws.Range("N108").Value = -97672.5
$ws.Range("H117").Value = 66063.5
$ws.Range("J117").Value = 66063.5
$ws.Range("L117").Value = 66063.5
$ws.Range("N117").Value = -75241.5
$ws.Range("H121").Value = 77300.75
$ws.Range("J121").Value = 77300.75
$ws.Range("L121").Value = 77300.75
$ws.Range("N121").Value = -80794.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 34995
$ws.Range("I2").Value = 34985
$ws.Range("J2").Value = 34997.5
$ws.Range("K2").Value = 34985
$ws.Range("L2").Value = 34997.5
$ws.Range("M2").Value = -34872
$ws.Range("N2").Value = -35223.5
$ws.Range("H52").Value = 77912
$ws.Range("J52").Value = 77912
$ws.Range("L52").Value = 77912
$ws.Range("N52").Value = -78438
$ws.Range("H53").Value = 38798.6
$ws.Range("J53").Value = 38821
$ws.Range("L53").Value = 38821
$ws.Range("N53").Value = -39969
$ws.Range("H86").Value = 2742.875
$ws.Range("I86").Value = 1957.8
$ws.Range("J86").Value = 3099.7273
$ws.Range("K86").Value = 1957.8
$ws.Range("L86").Value = 3099.7273
$ws.Range("M86").Value = -834.8
$ws.Range("N86").Value = -5345.7273
$ws.Range("H89").Value = 2742.875
$ws.Range("I89").Value = 1957.8
$ws.Range("J89").Value = 3099.7273
$ws.Range("K89").Value = 9789
$ws.Range("L89").Value = 15498.6365
$ws.Range("M89").Value = -4173
$ws.Range("N89").Value = -26730.6365
$ws.Range("H99").Value = 1937.5
$ws.Range("I99").Value = 1661.1765
$ws.Range("J99").Value = 3503.3333
$ws.Range("K99").Value = 1661.1765
$ws.Range("L99").Value = 3503.3333
$ws.Range("M99").Value = -163.1765
$ws.Range("N99").Value = -6499.3333
$ws.Range("H107").Value = 2060.2917
$ws.Range("I107").Value = 1711.3846
$ws.Range("K107").Value = 1711.3846
$ws.Range("M107").Value = 208.6153999999999
$ws.Range("H117").Value = 82941.8
$ws.Range("J117").Value = 82941.8
$ws.Range("L117").Value = 82941.8
$ws.Range("N117").Value = -92119.8
$ws.Range("H121").Value = 77912
$ws.Range("J121").Value = 77912
$ws.Range("L121").Value = 77912
$ws.Range("N121").Value = -81406
$ws.Range("H122").Value = 77996
$ws.Range("J122").Value = 77996
$ws.Range("L122").Value = 77996
$ws.Range("N122").Value = -87796
$ws.Range("H132").Value = 29200
$ws.Range("J132").Value = 29200
$ws.Range("L132").Value = 29200
$ws.Range("N132").Value = -39320

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4516.5317
$ws.Range("I31").Value = 2539.423
$ws.Range("J31").Value = 6964.381
$ws.Range("K31").Value = 2539.423
$ws.Range("L31").Value = 6964.381
$ws.Range("M31").Value = -2244.423
$ws.Range("N31").Value = -7554.381
$ws.Range("H34").Value = 4516.5317
$ws.Range("I34").Value = 2539.423
$ws.Range("J34").Value = 6964.381
$ws.Range("K34").Value = 2539.423
$ws.Range("L34").Value = 6964.381
$ws.Range("M34").Value = -2337.423
$ws.Range("N34").Value = -7368.381
$ws.Range("H99").Value = 2663.8572
$ws.Range("I99").Value = 999
$ws.Range("J99").Value = 2941.3333
$ws.Range("K99").Value = 999
$ws.Range("L99").Value = 2941.3333
$ws.Range("M99").Value = 499
$ws.Range("N99").Value = -5937.3333
$ws.Range("H105").Value = 1595.7333
$ws.Range("I105").Value = 1356.6923
$ws.Range("K105").Value = 1356.6923
$ws.Range("M105").Value = 390.3077000000001
$ws.Range("H107").Value = 691.3684
$ws.Range("J107").Value = 574.125
$ws.Range("L107").Value = 574.125
$ws.Range("N107").Value = -4414.125
$ws.Range("H116").Value = 41139.215
$ws.Range("J116").Value = 41139.215
$ws.Range("L116").Value = 41139.215
$ws.Range("N116").Value = -50317.215
$ws.Range("H119").Value = 96991.664
$ws.Range("J119").Value = 96991.664
$ws.Range("L119").Value = 96991.664
$ws.Range("N119").Value = -106667.664
$ws.Range("H126").Value = 2663.8572
$ws.Range("I126").Value = 999
$ws.Range("J126").Value = 2941.3333
$ws.Range("K126").Value = 2997
$ws.Range("L126").Value = 8823.999899999999
$ws.Range("M126").Value = -527
$ws.Range("N126").Value = -13763.9999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 23.333334
$ws.Range("H38").Value = 42.136364
$ws.Range("I38").Value = 26.083334
$ws.Range("J38").Value = 61.4
$ws.Range("K38").Value = 78.25000199999999
$ws.Range("L38").Value = 184.2
$ws.Range("M38").Value = 268.749998
$ws.Range("N38").Value = -878.2
$ws.Range("H122").Value = 1419.9
$ws.Range("I122").Value = 979.8
$ws.Range("K122").Value = 8818.199999999999
$ws.Range("M122").Value = -6368.199999999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 931.2308
$ws.Range("I107").Value = 942.2
$ws.Range("J107").Value = 894.6667
$ws.Range("K107").Value = 942.2
$ws.Range("L107").Value = 894.6667
$ws.Range("M107").Value = 977.8
$ws.Range("N107").Value = -4734.6667
$ws.Range("H108").Value = 64936.668
$ws.Range("J108").Value = 64936.668
$ws.Range("L108").Value = 64936.668
$ws.Range("N108").Value = -72616.66800000001
$ws.Range("H114").Value = 98317.664
$ws.Range("J114").Value = 98317.664
$ws.Range("L114").Value = 98317.664
$ws.Range("N114").Value = -106995.664

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1926.6471
$ws.Range("I16").Value = 2193.3845
$ws.Range("J16").Value = 1059.75
$ws.Range("K16").Value = 2193.3845
$ws.Range("L16").Value = 1059.75
$ws.Range("M16").Value = -2023.3845
$ws.Range("N16").Value = -1399.75
$ws.Range("H97").Value = 12353.25
$ws.Range("J97").Value = 12353.25
$ws.Range("L97").Value = 12353.25
$ws.Range("N97").Value = -14335.25
$ws.Range("H116").Value = 205935.8
$ws.Range("J116").Value = 205935.8
$ws.Range("L116").Value = 205935.8
$ws.Range("N116").Value = -215113.8
$ws.Range("H123").Value = 84994.44500000001
$ws.Range("J123").Value = 84994.44500000001
$ws.Range("L123").Value = 84994.44500000001
$ws.Range("N123").Value = -94794.44500000001
$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000
$ws.Range("H136").Value = 5124.9707
$ws.Range("I136").Value = 5051.625
$ws.Range("K136").Value = 15154.875
$ws.Range("M136").Value = -12604.875

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 44499.5
$ws.Range("J44").Value = 44499.5
$ws.Range("L44").Value = 44499.5
$ws.Range("N44").Value = -45607.5
$ws.Range("H59").Value = 36991
$ws.Range("J59").Value = 36991
$ws.Range("L59").Value = 36991
$ws.Range("N59").Value = -38467
